# Weekly fruit/vegetable price update:
# A brand-new price record (dated 2021-11-04, serial 44504) is inserted as the
# newest entry for this market/category, pushing the existing historical rows
# (170-201) down by one so the data keeps its usual newest-first ordering,
# and the former last row (201) survives as the new last row (202).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 170..201 down to 171..202, inheriting formatting from row 170.
$ws.Rows.Item(170).Insert()

# The freshly-inserted row 170 is blank; seed it with the same record as the
# row right below it (the old row 170, now shifted to 171) so every column
# except the date carries forward unchanged.
$ws.Range("A171:R171").Copy()
$ws.Range("A170").PasteSpecial()

# Stamp the new record with its own (later) date.
$ws.Range("D170").Value = 44504
